$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# A1 keeps its text label (shared string), unchanged in content.
$ws.Range("A1").Value = "HK_R_acc_SD"

# Updated mean values for A2:A50 (one extra row appended at the bottom, row 50).
$values = @(
    96.24183006535948,
    96.40522875816994,
    96.514161220043576,
    96.023965141612194,
    96.623093681917211,
    96.895424836601308,
    95.588235294117652,
    95.642701525054463,
    96.023965141612194,
    96.078431372549019,
    92.755991285403056,
    97.167755991285404,
    95.915032679738559,
    96.296296296296291,
    95.697167755991288,
    96.023965141612194,
    96.677559912854022,
    96.732026143790847,
    97.113289760348593,
    96.895424836601308,
    96.677559912854022,
    96.078431372549019,
    96.078431372549019,
    95.751633986928113,
    96.514161220043576,
    96.350762527233115,
    96.350762527233115,
    96.786492374727672,
    96.514161220043576,
    96.296296296296291,
    96.568627450980387,
    96.895424836601308,
    97.004357298474943,
    96.078431372549019,
    96.13289760348583,
    90.305010893246191,
    96.023965141612194,
    96.459694989106765,
    96.13289760348583,
    96.187363834422655,
    96.40522875816994,
    96.40522875816994,
    96.023965141612194,
    96.078431372549019,
    96.514161220043576,
    96.459694989106765,
    94.880174291938999,
    96.568627450980387,
    85.110227683411637
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
